$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in TODO row 4 with a new priority number and task description.
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = "サーバのトップページにあるメニューやID/PASSWORD情報を消す。"

# The row used to have D4:E4 merged (with C4 left standalone/empty).
# Now that C4 holds the task text, merge C4:E4 like the other filled-in rows.
$ws.Range("D4:E4").UnMerge()
$ws.Range("C4:E4").Merge()

# Update the view: scroll back to the top and select B5 (instead of row 10).
$null = $ws.Range("B5").Select()
